$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the reminder times in column A (time-of-day serial values).
# Row 2: 10:14 AM -> 6:45 AM
$ws.Range("A2").Value = 0.28125
# Row 3: 10:15 AM -> 6:47 AM
$ws.Range("A3").Value = 0.28263888888888888
